# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# The "Periodo Mora" rows (16-22) are refreshed with the new database
# extract: periods now run newest-first (2110 down to 2104), and the
# "Valor Mora" amount that used to belong to period 2110 now travels
# with it to the top row, while the remaining periods keep the 36341
# amount.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @("2110", "2109", "2108", "2107", "2106", "2105", "2104")
$valores = @(30284, 36341, 36341, 36341, 36341, 36341, 36341)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
